# Swap the distinguishing data between row 26 and row 27
# (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture row 26 original values (Value2 avoids the Variant-descriptor bug in Value) ---
$a26 = $ws.Range("A26").Value2
$b26 = $ws.Range("B26").Value2
$e26 = $ws.Range("E26").Value2
$f26 = $ws.Range("F26").Value2
$g26 = $ws.Range("G26").Value2
$h26 = $ws.Range("H26").Value2
$q26 = $ws.Range("Q26").Value2
$r26 = $ws.Range("R26").Value2

# --- capture row 27 original values ---
$a27 = $ws.Range("A27").Value2
$b27 = $ws.Range("B27").Value2
$e27 = $ws.Range("E27").Value2
$f27 = $ws.Range("F27").Value2
$g27 = $ws.Range("G27").Value2
$h27 = $ws.Range("H27").Value2
$q27 = $ws.Range("Q27").Value2
$r27 = $ws.Range("R27").Value2

# --- write row 27's original values into row 26 ---
$ws.Range("A26").Value2 = $a27
$ws.Range("B26").Value2 = $b27
$ws.Range("E26").Value2 = $e27
$ws.Range("F26").Value2 = $f27
$ws.Range("G26").Value2 = $g27
$ws.Range("H26").Value2 = $h27
$ws.Range("Q26").Value2 = $q27
$ws.Range("R26").Value2 = $r27

# --- write row 26's original values into row 27 ---
$ws.Range("A27").Value2 = $a26
$ws.Range("B27").Value2 = $b26
$ws.Range("E27").Value2 = $e26
$ws.Range("F27").Value2 = $f26
$ws.Range("G27").Value2 = $g26
$ws.Range("H27").Value2 = $h26
$ws.Range("Q27").Value2 = $q26
$ws.Range("R27").Value2 = $r26
